{"js": "// Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n// \"\u00a9 2020 ...\" footer paragraph that used to follow the\n// \"LOQ4038: Qu\u00edmica Org\u00e2nica II (Requisito fraco)\" paragraph, while leaving\n// the paragraph that originally came right before the page-break paragraph.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\n// Find the index of the \"Ver no Jupiter...\" paragraph (the anchor for the\n// whole block being removed), then delete it, the \"\u00a9 2020 ...\" paragraph\n// right after it, and the blank paragraph right before it.\nlet jupiterIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === targets[0]) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  const toDelete = [];\n  // Blank paragraph immediately before \"Ver no Jupiter...\"\n  if (jupiterIndex - 1 >= 0 && paras.items[jupiterIndex - 1].text === \"\") {\n    toDelete.push(paras.items[jupiterIndex - 1]);\n  }\n  toDelete.push(paras.items[jupiterIndex]);\n  // The copyright paragraph immediately after \"Ver no Jupiter...\"\n  if (jupiterIndex + 1 < paras.items.length && paras.items[jupiterIndex + 1].text === targets[1]) {\n    toDelete.push(paras.items[jupiterIndex + 1]);\n  }\n  // Delete from the last one to the first one so earlier indices stay valid.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the empty paragraph, the \"Ver no Jupiter...\" paragraph, and the\n# \"\u00a9 2020 ...\" footer paragraph that used to follow the\n# \"LOQ4038: Qu\u00edmica Org\u00e2nica II (Requisito fraco)\" paragraph, while leaving\n# the paragraph that originally came right before the page-break paragraph.\n$d = $word.ActiveDocument\n\n$target = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyright = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$idx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd()\n    if ($t -eq $target) {\n        $idx = $i\n        break\n    }\n}\n\nif ($idx -ne -1) {\n    $beforeIdx = $idx - 1\n    $afterIdx = $idx + 1\n    $beforeText = \"\"\n    if ($beforeIdx -ge 1) {\n        $beforeText = $d.Paragraphs.Item($beforeIdx).Range.Text.TrimEnd()\n    }\n    $afterText = \"\"\n    if ($afterIdx -le $d.Paragraphs.Count) {\n        $afterText = $d.Paragraphs.Item($afterIdx).Range.Text.TrimEnd()\n    }\n\n    # Delete from the highest index down to the lowest so earlier indices\n    # stay valid while we work.\n    if ($afterText -eq $copyright) {\n        $d.Paragraphs.Item($afterIdx).Range.Delete()\n    }\n    $d.Paragraphs.Item($idx).Range.Delete()\n    if ($beforeText -eq \"\") {\n        $d.Paragraphs.Item($beforeIdx).Range.Delete()\n    }\n}\n"}
